# Edit: Added new TCs for Duplicate SO number, Avalara Tax retained TC...
# Target sheet: "ARATO" - insert a new leading column (CompanyID) that was
# missing, shifting the existing 19 columns (A:S) one place to the right
# (B:T), then populate the new column A with the CompanyID values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARATO")

# Insert a brand-new column before column A; everything currently in
# A:S shifts right into B:T.
$ws.Columns.Item(1).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 1).Value2 = "CompanyID"

# Populate the CompanyID value for every data row (2-17). This mirrors
# the value already present in what is now column E (FinancialCompany /
# "aBb5f0000004JfX") for each row.
for ($r = 2; $r -le 17; $r++) {
    $v = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 1).Value2 = $v
}

$ws.Range("B23").Select()
